$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.562.12"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.612.63"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'511.43"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'154.68"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "2.630.97"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "3.073.52"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "60.517.97"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "2.628.39"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "'353.13"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'60.68"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  -3.30%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'151.55"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").Value = "'0.891"
$ws.Range("E37").Value = "  +4.99%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "'291.12"
$ws.Range("E42").Value = "  -5.99%  "
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").Value = "'19.86"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "'4.92"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "'10.31"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "1.998.56"
$ws.Range("E51").Value = "  -3.17%  "

# Row 39/40 swap: OKB <-> Fetch.AI with new data
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'36.35"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'0.844"
$ws.Range("E40").Value = "  -2.00%  "
